$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 811.6
$ws.Cells.Item(17, 10).Value = 811.6
$ws.Cells.Item(17, 12).Value = 2434.8
$ws.Cells.Item(17, 14).Value = -2770.8
$ws.Cells.Item(32, 8).Value = 7666.952
$ws.Cells.Item(32, 10).Value = 5034.3335
$ws.Cells.Item(32, 12).Value = 5034.3335
$ws.Cells.Item(32, 14).Value = -5686.3335
$ws.Cells.Item(51, 8).Value = 6725.25
$ws.Cells.Item(51, 10).Value = 8061.4
$ws.Cells.Item(51, 12).Value = 8061.4
$ws.Cells.Item(51, 14).Value = -9029.4
$ws.Cells.Item(55, 8).Value = 14919.4
$ws.Cells.Item(55, 9).Value = 531.6667
$ws.Cells.Item(55, 10).Value = 36501
$ws.Cells.Item(55, 11).Value = 531.6667
$ws.Cells.Item(55, 12).Value = 36501
$ws.Cells.Item(55, 13).Value = -317.6667
$ws.Cells.Item(55, 14).Value = -36929
$ws.Cells.Item(137, 8).Value = 15769.454
$ws.Cells.Item(137, 9).Value = 20470.812
$ws.Cells.Item(137, 10).Value = 3232.5
$ws.Cells.Item(137, 11).Value = 61412.436
$ws.Cells.Item(137, 12).Value = 9697.5
$ws.Cells.Item(137, 13).Value = -58862.436
$ws.Cells.Item(137, 14).Value = -14797.5

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(4, 8).Value = 300
$ws.Cells.Item(4, 9).Value = 300
$ws.Cells.Item(4, 11).Value = 300
$ws.Cells.Item(4, 13).Value = -184
$ws.Cells.Item(5, 8).Value = 7064.467
$ws.Cells.Item(5, 9).Value = 8768.5
$ws.Cells.Item(5, 10).Value = 248.33333
$ws.Cells.Item(5, 11).Value = 8768.5
$ws.Cells.Item(5, 12).Value = 248.33333
$ws.Cells.Item(5, 13).Value = -8656.5
$ws.Cells.Item(5, 14).Value = -472.33333
$ws.Cells.Item(32, 8).Value = 34763.47
$ws.Cells.Item(32, 9).Value = 34763.47
$ws.Cells.Item(32, 11).Value = 34763.47
$ws.Cells.Item(32, 13).Value = -34476.47
$ws.Cells.Item(45, 8).Value = 3868.2
$ws.Cells.Item(45, 9).Value = 2505.4285
$ws.Cells.Item(45, 10).Value = 5060.625
$ws.Cells.Item(45, 11).Value = 2505.4285
$ws.Cells.Item(45, 12).Value = 5060.625
$ws.Cells.Item(45, 13).Value = -2128.4285
$ws.Cells.Item(45, 14).Value = -5814.625
$ws.Cells.Item(110, 8).Value = 1169.4286
$ws.Cells.Item(110, 10).Value = 1499
$ws.Cells.Item(110, 12).Value = 1499
$ws.Cells.Item(110, 14).Value = -5589
$ws.Cells.Item(122, 8).Value = 1800.7826
$ws.Cells.Item(122, 9).Value = 1705.9474
$ws.Cells.Item(122, 10).Value = 2251.25
$ws.Cells.Item(122, 11).Value = 5117.8422
$ws.Cells.Item(122, 12).Value = 6753.75
$ws.Cells.Item(122, 13).Value = -2667.8422
$ws.Cells.Item(122, 14).Value = -11653.75
$ws.Cells.Item(132, 8).Value = 1468.9452
$ws.Cells.Item(132, 9).Value = 1083.4894
$ws.Cells.Item(132, 10).Value = 2165.7307
$ws.Cells.Item(132, 11).Value = 3250.4682
$ws.Cells.Item(132, 12).Value = 6497.1921
$ws.Cells.Item(132, 13).Value = -720.4681999999998
$ws.Cells.Item(132, 14).Value = -11557.1921

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(4, 8).Value = 7064.467
$ws.Cells.Item(4, 9).Value = 8768.5
$ws.Cells.Item(4, 10).Value = 248.33333
$ws.Cells.Item(4, 11).Value = 8768.5
$ws.Cells.Item(4, 12).Value = 248.33333
$ws.Cells.Item(4, 13).Value = -8653.5
$ws.Cells.Item(4, 14).Value = -478.33333
$ws.Cells.Item(20, 8).Value = 15451.952
$ws.Cells.Item(20, 9).Value = 21263.467
$ws.Cells.Item(20, 11).Value = 21263.467
$ws.Cells.Item(20, 13).Value = -21016.467
$ws.Cells.Item(86, 8).Value = 4478.8335
$ws.Cells.Item(86, 9).Value = 1749.5
$ws.Cells.Item(86, 11).Value = 1749.5
$ws.Cells.Item(86, 13).Value = -626.5
$ws.Cells.Item(89, 8).Value = 4478.8335
$ws.Cells.Item(89, 9).Value = 1749.5
$ws.Cells.Item(89, 11).Value = 9000
$ws.Cells.Item(89, 13).Value = -3131.5
$ws.Cells.Item(120, 8).Value = 75979
$ws.Cells.Item(120, 10).Value = 75979
$ws.Cells.Item(120, 12).Value = 75979
$ws.Cells.Item(120, 14).Value = -85655
$ws.Cells.Item(134, 8).Value = 1545.7885
$ws.Cells.Item(134, 9).Value = 1112.8864
$ws.Cells.Item(134, 11).Value = 3338.6592
$ws.Cells.Item(134, 13).Value = -803.6592000000001

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 2153.158
$ws.Cells.Item(16, 9).Value = 1509.8572
$ws.Cells.Item(16, 11).Value = 1509.8572
$ws.Cells.Item(16, 13).Value = -1222.8572
$ws.Cells.Item(31, 8).Value = 4547658.5
$ws.Cells.Item(31, 9).Value = 6667665
$ws.Cells.Item(31, 10).Value = 4786.857
$ws.Cells.Item(31, 11).Value = 6667665
$ws.Cells.Item(31, 12).Value = 4786.857
$ws.Cells.Item(31, 13).Value = -6667370
$ws.Cells.Item(31, 14).Value = -5376.857
$ws.Cells.Item(34, 8).Value = 4547658.5
$ws.Cells.Item(34, 9).Value = 6667665
$ws.Cells.Item(34, 10).Value = 4786.857
$ws.Cells.Item(34, 11).Value = 6667665
$ws.Cells.Item(34, 12).Value = 4786.857
$ws.Cells.Item(34, 13).Value = -6667463
$ws.Cells.Item(34, 14).Value = -5190.857
$ws.Cells.Item(86, 8).Value = 52993.53
$ws.Cells.Item(86, 9).Value = 64991
$ws.Cells.Item(86, 11).Value = 64991
$ws.Cells.Item(86, 13).Value = -63868
$ws.Cells.Item(89, 8).Value = 52993.53
$ws.Cells.Item(89, 9).Value = 64991
$ws.Cells.Item(89, 11).Value = 324955
$ws.Cells.Item(89, 13).Value = -319339
$ws.Cells.Item(113, 8).Value = 2153.158
$ws.Cells.Item(113, 9).Value = 1509.8572
$ws.Cells.Item(113, 11).Value = 1509.8572
$ws.Cells.Item(113, 13).Value = 660.1428000000001

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(37, 8).Value = 42369.066
$ws.Cells.Item(37, 10).Value = 42369.066
$ws.Cells.Item(37, 12).Value = 127107.198
$ws.Cells.Item(37, 14).Value = -127331.198
$ws.Cells.Item(68, 8).Value = 4581.758
$ws.Cells.Item(68, 10).Value = 4790.3228
$ws.Cells.Item(68, 12).Value = 14370.9684
$ws.Cells.Item(68, 14).Value = -15992.9684
$ws.Cells.Item(71, 8).Value = 4581.758
$ws.Cells.Item(71, 10).Value = 4790.3228
$ws.Cells.Item(71, 12).Value = 43112.9052
$ws.Cells.Item(71, 14).Value = -51224.9052
$ws.Cells.Item(133, 8).Value = 3109.3809
$ws.Cells.Item(133, 9).Value = 2296.1875
$ws.Cells.Item(133, 10).Value = 5711.6
$ws.Cells.Item(133, 11).Value = 6888.5625
$ws.Cells.Item(133, 12).Value = 17134.8
$ws.Cells.Item(133, 13).Value = -1828.5625
$ws.Cells.Item(133, 14).Value = -27254.8

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(32, 8).Value = 27496.334
$ws.Cells.Item(32, 10).Value = 27496.334
$ws.Cells.Item(32, 12).Value = 27496.334
$ws.Cells.Item(32, 14).Value = -28088.334
$ws.Cells.Item(97, 8).Value = 621.63635
$ws.Cells.Item(97, 9).Value = 747.73334
$ws.Cells.Item(97, 10).Value = 351.42856
$ws.Cells.Item(97, 11).Value = 747.73334
$ws.Cells.Item(97, 12).Value = 351.42856
$ws.Cells.Item(97, 13).Value = -251.73334
$ws.Cells.Item(97, 14).Value = -1343.42856
$ws.Cells.Item(107, 8).Value = 417.625
$ws.Cells.Item(107, 9).Value = 122.25
$ws.Cells.Item(107, 10).Value = 713
$ws.Cells.Item(107, 11).Value = 122.25
$ws.Cells.Item(107, 12).Value = 713
$ws.Cells.Item(107, 13).Value = 1797.75
$ws.Cells.Item(107, 14).Value = -4553
$ws.Cells.Item(132, 8).Value = 2031.0883
$ws.Cells.Item(132, 9).Value = 1482.1154
$ws.Cells.Item(132, 11).Value = 4446.3462
$ws.Cells.Item(132, 13).Value = -1916.3462

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(11, 8).Value = 28833.334
$ws.Cells.Item(11, 10).Value = 28833.334
$ws.Cells.Item(11, 12).Value = 28833.334
$ws.Cells.Item(11, 14).Value = -29113.334
$ws.Cells.Item(82, 8).Value = 1552.0834
$ws.Cells.Item(82, 9).Value = 1296.8572
$ws.Cells.Item(82, 11).Value = 1296.8572
$ws.Cells.Item(82, 13).Value = -935.8571999999999
$ws.Cells.Item(85, 8).Value = 1552.0834
$ws.Cells.Item(85, 9).Value = 1296.8572
$ws.Cells.Item(85, 11).Value = 1296.8572
$ws.Cells.Item(85, 13).Value = -48.85719999999992

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(64, 8).Value = 100000
$ws.Cells.Item(64, 10).Value = 100000
$ws.Cells.Item(64, 12).Value = 100000
$ws.Cells.Item(64, 14).Value = -100496
$ws.Cells.Item(67, 8).Value = 100000
$ws.Cells.Item(67, 10).Value = 100000
$ws.Cells.Item(67, 12).Value = 100000
$ws.Cells.Item(67, 14).Value = -101716
$ws.Cells.Item(93, 8).Value = 100000
$ws.Cells.Item(93, 10).Value = 100000
$ws.Cells.Item(93, 12).Value = 100000
$ws.Cells.Item(93, 14).Value = -104992
$ws.Cells.Item(122, 8).Value = 52228.38
$ws.Cells.Item(122, 9).Value = 64680.305
$ws.Cells.Item(122, 11).Value = 194040.915
$ws.Cells.Item(122, 13).Value = -191590.915
$ws.Cells.Item(132, 8).Value = 1735.2963
$ws.Cells.Item(132, 9).Value = 1292.75
$ws.Cells.Item(132, 11).Value = 3878.25
$ws.Cells.Item(132, 13).Value = -1348.25
